$d = $word.ActiveDocument

# Locate the paragraph ending in "Enjoy the story!" (the about-page closer) so the
# new "lore dump" paragraphs can be inserted immediately after it, regardless of
# its exact paragraph index.
$findRange = $d.Content
$ok = $findRange.Find.Execute("Enjoy the story!", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $ok) {
    throw "Could not find anchor paragraph ending in `"Enjoy the story!`""
}
$findRange.Collapse(0)
$idx = $findRange.Paragraphs.Item(1).Index

# New paragraph texts to insert, in order (empty string = blank paragraph).
$newParagraphs = @(
    ""
    "LORE SO FAR"
    "NOTE: The following page contains spoilers."
    "Because this story is still in development, only a very small amount of lore has been established so far. There are a series of survival games (similar to Saw or Squid Game) that the players in this universe are subjected to. Additional hazards and traps like the blistering heat, scorpions, poisons, bombs, and land mines are scattered throughout the wasteland. Resources like food, water, first aid, shelter, and tools are deliberately scarce. Some players may be well equipped with items and acutely aware of a particular game’s rules including how to win, or be at a complete disadvantage armed with nothing and given no information at all. Not much is known about who or how the game has been orchestrated, but it is generally advised not to interfere or collaborate with other players unless instructed to… "
    "The first game: the player must kill their target before their opponent does or the collar around their neck explodes. The target player assigned matches the number tattooed over the left eye. Once the player kills their target, their collar will release and their opponent’s will detonate. "
    "Other than that, that’s all I’ve got."
)

foreach ($text in $newParagraphs) {
    $p = $d.Paragraphs.Item($idx)
    $p.Range.InsertParagraphAfter()
    $idx = $idx + 1
    if ($text -ne "") {
        $d.Paragraphs.Item($idx).Range.Text = $text
    }
}
